# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scraped totals (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Row => New Value map (F column) shared by both sheets
$updates = @{
    2  = 294
    4  = 864
    5  = 22
    6  = 315
    7  = 9416
    8  = 79
    9  = 80
    10 = 132
    11 = 120
    12 = 11
    14 = 26
    17 = 276
    18 = 776
    19 = 44
    20 = 93
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
